$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New issue #21 (FR_HIRING / enhancement) goes into row 17.
# Column A ("21") looks numeric, so a plain assignment would store it as a
# number; prefix with an apostrophe to force text entry (matching the other
# rows, which all store the Issue ID as text), then clear the resulting
# "quote prefix" cell format so no extra style is introduced.
$ws.Range("A17").Value = "'21"
$ws.Range("A17").ClearFormats()

$ws.Range("B17").Value = "FR_HIRING"
$ws.Range("C17").Value = "open"
$ws.Range("D17").Value = "2025-03-26T06:41:40Z"
$ws.Range("E17").Value = "enhancement"
